$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.029.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +2.05%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.672.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +2.91%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'216.11"
$ws.Range("D5").Style = "Normal"

$ws.Range("E6").Value = "'  +2.07%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = "'  -0.02%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +2.05%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D10").Value = "'20.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +4.87%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.0893"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +4.66%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.908.26"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +2.90%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'1.667.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +2.43%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'4.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.26%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("E16").Value = "'  +1.75%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'27.043.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +2.09%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'234.90"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.01%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'0.0₃0736"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +1.33%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  -0.80%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("E21").Value = "'  +0.11%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'4.46"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +3.40%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'2.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.18%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  +1.22%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'145.25"
$ws.Range("D25").Style = "Normal"

$ws.Range("E26").Value = "'  +1.33%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("E27").Value = "'  +0.55%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'15.89"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +1.49%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -0.07%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = "'  -0.02%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("E31").Value = "'  +1.45%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = "'  +2.17%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.452.03"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -4.67%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'3.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +5.18%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = "'  +5.40%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = "'  -0.30%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").Value = "'ARBITRUM"
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").Value = "'0.893"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +6.90%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").Value = "'ImmutableX"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").Value = "'0.569"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.27%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.0170"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +1.86%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +3.72%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("E41").Value = "'  +0.02%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("E42").Value = "'  +11.61%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = "'  +2.70%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'65.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.27%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'1.816.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.90%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.781"
$ws.Range("D46").Style = "Normal"

$ws.Range("D47").Value = "'90.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +0.21%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  +1.37%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("E49").Value = "'  +4.15%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  +1.56%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'7.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.50%  "
$ws.Range("E51").Style = "Normal"
